$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.134.03"
$ws.Range("E2").Value = "  -0.66%  "

# Row 3
$ws.Range("D3").Value = "2.469.21"
$ws.Range("E3").Value = "  -2.34%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "583.20"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.52%  "

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "169.05"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.58%  "

# Row 7
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.514"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -2.14%  "

# Row 9
$ws.Range("D9").Value = "2.470.98"
$ws.Range("E9").Value = "  -2.21%  "

# Row 10
$ws.Range("E10").Value = "  -2.53%  "

# Row 11
$ws.Range("E11").Value = "  -0.01%  "

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "4.97"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -2.74%  "

# Row 13
$ws.Range("E13").Value = "  -3.33%  "

# Row 14
$ws.Range("E14").Value = "  -3.04%  "

# Row 15
$ws.Range("D15").Value = "2.866.45"
$ws.Range("E15").Value = "  -2.07%  "

# Row 16
$ws.Range("D16").Value = "67.053.81"
$ws.Range("E16").Value = "  -0.64%  "

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.0000169"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -4.19%  "

# Row 18
$ws.Range("D18").Value = "2.465.29"
$ws.Range("E18").Value = "  -1.75%  "

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.16"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -5.36%  "

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.62"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -3.30%  "

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "353.86"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -4.00%  "

# Row 22
$ws.Range("E22").Value = "  -2.58%  "

# Row 23
$ws.Range("E23").Value = "  -0.01%  "

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "69.06"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -3.61%  "

# Row 25
$ws.Range("E25").Value = "  -7.27%  "

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "1.78"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -7.08%  "

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "9.28"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -6.76%  "

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.997"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.31%  "

# Row 29
$ws.Range("D29").Value = "2.588.70"
$ws.Range("E29").Value = "  -2.47%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0907"
$ws.Range("E30").Value = "  -5.91%  "

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "516.42"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -3.89%  "

# Row 32
$ws.Range("E32").Value = "  -7.23%  "

# Row 33
$ws.Range("E33").Value = "  -5.18%  "

# Row 34
$ws.Range("E34").Value = "  -5.54%  "

# Row 35
$ws.Range("E35").Value = "  +0.07%  "

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.120"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -6.66%  "

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "157.68"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.57%  "

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "18.68"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.29%  "

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "18.40"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -3.27%  "

# Row 40
$ws.Range("E40").Value = "  -5.30%  "

# Row 41
$ws.Range("E41").Value = "  +0.24%  "

# Row 42
$ws.Range("E42").Value = "  -6.80%  "

# Row 43
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.66"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -6.33%  "

# Row 44
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "4.79"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -6.50%  "

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "2.37"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -5.58%  "

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "38.73"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -2.20%  "

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "141.16"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.56%  "

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "3.46"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -6.59%  "

# Row 49
$ws.Range("E49").Value = "  -6.44%  "

# Row 50
$ws.Range("E50").Value = "  -11.64%  "

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.59"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -7.17%  "

Write-Host "Update complete"